$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.962.69'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").Value = '2.359.85'
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E5").Value = '  +0.85%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '240.19'
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '74.61'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +2.80%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +11.02%  '
$ws.Range("E10").Value = '  +0.85%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '57.22'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.07%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '32.31'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +11.29%  '
$ws.Range("E13").Value = '  +9.85%  '
$ws.Range("E14").Value = '  +0.80%  '
$ws.Range("D15").Value = '2.709.69'
$ws.Range("E15").Value = '  +0.33%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '16.66'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.66%  '
$ws.Range("E17").Value = '  -0.02%  '
$ws.Range("D18").Value = '2.361.63'
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("D19").Value = '43.923.64'
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("E21").Value = '  +5.44%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '77.11'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -1.10%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '258.33'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +1.57%  '
$ws.Range("E24").Value = '  +25.99%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.51'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.09%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '3.67'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -2.21%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '10.80'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +3.11%  '
$ws.Range("E29").Value = '  -0.21%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '22.79'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.86%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '175.73'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.81%  '
$ws.Range("E32").Value = '  -1.88%  '
$ws.Range("E33").Value = '  +3.48%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.0769'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +7.08%  '
$ws.Range("E35").Value = '  +1.76%  '
$ws.Range("E36").Value = '  +4.86%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '3.80'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -3.66%  '
$ws.Range("E38").Value = '  -2.51%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '6.34'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.23%  '
$ws.Range("E40").Value = '  +5.25%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.113'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +15.21%  '
$ws.Range("E42").Value = '  +15.09%  '
$ws.Range("E43").Value = '  +3.54%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '19.14'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -1.46%  '
$ws.Range("E45").Value = '  +0.06%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '4.75'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +6.77%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '2.54'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +9.13%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '58.21'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +11.27%  '
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("E50").Value = '  +1.20%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '100.48'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +2.30%  '
